# Generate Report for Handoff
# The f9f150d6-... file has moved from "In Translation" to "Ready for
# handoff" with a new handoff timestamp. Update the Overview sheet and the
# per-locale (zh-cn / de-de) sheets accordingly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet, row 3 (f9f150d6 file) -------------------------------
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 20:11:56"

# ---- zh-cn sheet, row 3 (f9f150d6 file) -----------------------------------
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-17 20:11:51"

# ---- de-de sheet, row 3 (f9f150d6 file) -----------------------------------
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-17 20:11:56"

# ---- Column width follow-up (status/datetime columns grew wider to fit
# ---- the new "Ready for handoff" text) ------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
